# "for function for adding tenants" -- add a 4th column ("add tenants") with a
# small list of tenant names / contact info, fix the B4 typo'd value, and give
# the whole table a consistent centered look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data fix: B4 was 210, should be 200
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 200

# ---------------------------------------------------------------------------
# 2. New column D -- header + 10 rows of tenant data, then a few blank
#    formatted rows to round out the new table to row 14.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "add tenants"

$ws.Range("D2").Value = "יוסי"
$ws.Range("D3").Value = "0524455586"
$ws.Range("D4").Value = "qa@oxs.co.il"
$ws.Range("D5").Value = "משה"
$ws.Range("D6").Value = "דני"
$ws.Range("D7").Value = "דוד"
$ws.Range("D8").Value = "רמי"
$ws.Range("D9").Value = "לאה"
$ws.Range("D10").Value = "שרה"
$ws.Range("D11").Value = "חוה"

# D12:D14 stay empty but still get formatted like the rest of the column.

# ---------------------------------------------------------------------------
# 3. Column width for the new column (matches the rest of the sheet's
#    explicit custom widths).
# ---------------------------------------------------------------------------
$ws.Range("D1").ColumnWidth = 17

# ---------------------------------------------------------------------------
# 4. Borders: give every cell in the (now 14-row) table a thin box border so
#    the new rows/column match the look of the existing data rows.
# ---------------------------------------------------------------------------
$body = $ws.Range("A3:D14")
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2

$d2 = $ws.Range("D2")
$d2.Borders.LineStyle = 1
$d2.Borders.Weight = 2

# ---------------------------------------------------------------------------
# 5. Hyperlink on D4 (the new "qa@oxs.co.il" contact), and make D7/D10 look
#    the same (underlined hyperlink-blue) even though they are plain text.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:qa@oxs.co.il")

$ws.Range("D7").Font.Underline = 2
$ws.Range("D7").Font.ThemeColor = 11
$ws.Range("D10").Font.Underline = 2
$ws.Range("D10").Font.ThemeColor = 11
$ws.Range("D13").Font.Underline = 2
$ws.Range("D13").Font.ThemeColor = 11

# ---------------------------------------------------------------------------
# 6. Number format: the phone-number-ish column D cells are stored as Text
#    so leading zeros / formatting are preserved (D3, D6, D9, D12).
# ---------------------------------------------------------------------------
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 7. Alignment: center (both horizontal + vertical) across the whole table,
#    header included.
# ---------------------------------------------------------------------------
$whole = $ws.Range("A1:D14")
$whole.HorizontalAlignment = -4108
$whole.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 8. Selection marker, matching where the author's cursor ended up.
# ---------------------------------------------------------------------------
$ws.Range("D12").Select()
